$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 481, pushing the existing rows 481:526 down to 483:528.
$ws.Rows("481:482").Insert()

# Row 481 (new week, Primera)
$ws.Cells.Item(481, 1).Value = 8
$ws.Cells.Item(481, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(481, 3).Value = "Coquimbo"
$ws.Cells.Item(481, 4).Value = 45223
$ws.Cells.Item(481, 5).Value = 4
$ws.Cells.Item(481, 6).Value = 100114014
$ws.Cells.Item(481, 7).Value = "Betarraga"
$ws.Cells.Item(481, 8).Value = "Sin especificar"
$ws.Cells.Item(481, 9).Value = "Primera"
$ws.Cells.Item(481, 10).Value = 1600
$ws.Cells.Item(481, 11).Value = 500
$ws.Cells.Item(481, 12).Value = 600
$ws.Cells.Item(481, 13).Value = 550
$ws.Cells.Item(481, 14).Value = '$/paquete 3 unidades'
$ws.Cells.Item(481, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(481, 16).Value = 183
$ws.Cells.Item(481, 17).Value = 3
$ws.Cells.Item(481, 18).Value = "Hortaliza"

# Row 482 (new week, Segunda)
$ws.Cells.Item(482, 1).Value = 8
$ws.Cells.Item(482, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(482, 3).Value = "Coquimbo"
$ws.Cells.Item(482, 4).Value = 45223
$ws.Cells.Item(482, 5).Value = 4
$ws.Cells.Item(482, 6).Value = 100114014
$ws.Cells.Item(482, 7).Value = "Betarraga"
$ws.Cells.Item(482, 8).Value = "Sin especificar"
$ws.Cells.Item(482, 9).Value = "Segunda"
$ws.Cells.Item(482, 10).Value = 1100
$ws.Cells.Item(482, 11).Value = 400
$ws.Cells.Item(482, 12).Value = 450
$ws.Cells.Item(482, 13).Value = 425
$ws.Cells.Item(482, 14).Value = '$/paquete 3 unidades'
$ws.Cells.Item(482, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(482, 16).Value = 142
$ws.Cells.Item(482, 17).Value = 3
$ws.Cells.Item(482, 18).Value = "Hortaliza"
